$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 6 into row 7 (copy keeps formatting, e.g. the date number
# format on column D), then overwrite row 6 with this week's new entry.
$ws.Rows("6:6").Copy()
$ws.Rows("7:7").Insert()

# Populate the new row 6 with this week's price entry.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Macroferia Regional de Talca"
$ws.Range("C6").Value = "Maule"
$ws.Range("D6").Value = 45275
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101004
$ws.Range("J6").Value = "Frambuesa"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 5000
$ws.Range("O6").Value = 5000
$ws.Range("P6").Value = 5000
$ws.Range("Q6").Value = "`$/bandeja 2 kilos"
$ws.Range("R6").Value = "Región del Maule"
$ws.Range("S6").Value = 2500
$ws.Range("T6").Value = 2
